{"js": "// Fix badly worded sentence: \"es una de tarea complicada\" -> \"es una tarea complicada\"\n// (removes the erroneous \"de\" before \"tarea\").\nconst searchText = \"es una de tarea complicada\";\nconst results = context.document.body.search(searchText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found: \" + searchText);\n}\n\nfor (const r of results.items) {\n  r.insertText(\"es una tarea complicada\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Arreglo algo que esta mal redactado\n# Fix badly worded sentence: remove the erroneous \"de\" so\n# \"es una de tarea complicada\" reads \"es una tarea complicada\".\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"es una de tarea complicada\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"es una tarea complicada\"\n\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    throw \"Target text not found: 'es una de tarea complicada'\"\n}\n"}
